# Fix the big mistake: correct Total (col B) and Community (col D) values
# for rows 2-13 of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = @{ B = 2993.3023958;      D = 202.2209602333333 }
    3  = @{ B = 2815.513484116667; D = 193.8054265166667 }
    4  = @{ B = 3002.775658366667; D = 200.7630757 }
    5  = @{ B = 2896.734895566667; D = 202.4454222333333 }
    6  = @{ B = 2986.844631983334; D = 205.8283510833333 }
    7  = @{ B = 2939.527072450001; D = 187.4918758 }
    8  = @{ B = 2981.484152850001; D = 195.3420448 }
    9  = @{ B = 3008.645255083334; D = 203.6077692666667 }
    10 = @{ B = 2929.975138683334; D = 195.5894932833333 }
    11 = @{ B = 3006.062622700001; D = 194.65602 }
    12 = @{ B = 2915.079398966667; D = 205.13571825 }
    13 = @{ B = 2931.312351233334; D = 189.0492176666667 }
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Range("B$row").Value = $rowValues.B
    $ws.Range("D$row").Value = $rowValues.D
}
